# Burundi final test - update workbook:
#  - Remove "BIGWA II" school row
#  - Insert 8 new school rows (MUTUMBA II, MUYANGE I, MUYANGE II, MUYUGA,
#    MUZIMA, MWAZA, MWUMBA, NDAGO) after BIKINGI / before BUBAJI I
#  - Resize the "Tabelle1" table/autofilter to the new extent
#  - Make "Schools" the active sheet/tab, with B21 selected

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schools")
$ws2 = $wb.Worksheets.Item("Warehouses")

# 1. Remove the "BIGWA II" row (original row 5)
$ws1.Rows.Item(5).Delete()

# 2. Make room for the 8 new rows right after BIKINGI (now row 5)
$ws1.Rows("6:13").Insert()

# 3. Fill in the new rows with the new school data
#    columns: A Name_ID | B Total Sum of Beneficiaries | C Total Sum of Commodities
#             D Consumption per day in mt | E Consumption per week in mt
#             F Latitude | G Longitude | H Capacity | I Lower | J Initial | K Storage Cost

# Row 6 - MUTUMBA II
$ws1.Cells.Item(6, 1).Value = "MUTUMBA II"
$ws1.Cells.Item(6, 2).Value = 635
$ws1.Cells.Item(6, 3).Value = 1.933575
$ws1.Cells.Item(6, 4).Value = 0.127
$ws1.Cells.Item(6, 5).Value = 0.63500000000000001
$ws1.Cells.Item(6, 6).Value = -3.5886499999999999
$ws1.Cells.Item(6, 7).Value = 29.35942
$ws1.Cells.Item(6, 8).Value = 3.81
$ws1.Cells.Item(6, 9).Value = 0.64
$ws1.Cells.Item(6, 10).Value = 1.76
$ws1.Cells.Item(6, 11).Value = 0

# Row 7 - MUYANGE I
$ws1.Cells.Item(7, 1).Value = "MUYANGE I"
$ws1.Cells.Item(7, 2).Value = 1168
$ws1.Cells.Item(7, 3).Value = 4.8400000000000007
$ws1.Cells.Item(7, 4).Value = 0.23360000000000003
$ws1.Cells.Item(7, 5).Value = 1.1680000000000001
$ws1.Cells.Item(7, 6).Value = -2.8280959999999999
$ws1.Cells.Item(7, 7).Value = 30.127786
$ws1.Cells.Item(7, 8).Value = 6.73
$ws1.Cells.Item(7, 9).Value = 1.17
$ws1.Cells.Item(7, 10).Value = 1.7
$ws1.Cells.Item(7, 11).Value = 0

# Row 8 - MUYANGE II
$ws1.Cells.Item(8, 1).Value = "MUYANGE II"
$ws1.Cells.Item(8, 2).Value = 590
$ws1.Cells.Item(8, 3).Value = 1.843
$ws1.Cells.Item(8, 4).Value = 0.11800000000000001
$ws1.Cells.Item(8, 5).Value = 0.59000000000000008
$ws1.Cells.Item(8, 6).Value = -3.1234920000000002
$ws1.Cells.Item(8, 7).Value = 29.34451
$ws1.Cells.Item(8, 8).Value = 2.46
$ws1.Cells.Item(8, 9).Value = 0.59
$ws1.Cells.Item(8, 10).Value = 1.33
$ws1.Cells.Item(8, 11).Value = 0

# Row 9 - MUYUGA
$ws1.Cells.Item(9, 1).Value = "MUYUGA"
$ws1.Cells.Item(9, 2).Value = 847
$ws1.Cells.Item(9, 3).Value = 1.2709999999999999
$ws1.Cells.Item(9, 4).Value = 0.16940000000000002
$ws1.Cells.Item(9, 5).Value = 0.84700000000000009
$ws1.Cells.Item(9, 6).Value = -3.7038009999999999
$ws1.Cells.Item(9, 7).Value = 29.834140999999999
$ws1.Cells.Item(9, 8).Value = 4.83
$ws1.Cells.Item(9, 9).Value = 0.85
$ws1.Cells.Item(9, 10).Value = 1.8
$ws1.Cells.Item(9, 11).Value = 0

# Row 10 - MUZIMA
$ws1.Cells.Item(10, 1).Value = "MUZIMA"
$ws1.Cells.Item(10, 2).Value = 848
$ws1.Cells.Item(10, 3).Value = 2.7759999999999998
$ws1.Cells.Item(10, 4).Value = 0.16960000000000003
$ws1.Cells.Item(10, 5).Value = 0.84800000000000009
$ws1.Cells.Item(10, 6).Value = -3.5700080000000001
$ws1.Cells.Item(10, 7).Value = 29.841771999999999
$ws1.Cells.Item(10, 8).Value = 4.47
$ws1.Cells.Item(10, 9).Value = 0.85
$ws1.Cells.Item(10, 10).Value = 1.1100000000000001
$ws1.Cells.Item(10, 11).Value = 0

# Row 11 - MWAZA
$ws1.Cells.Item(11, 1).Value = "MWAZA"
$ws1.Cells.Item(11, 2).Value = 708
$ws1.Cells.Item(11, 3).Value = 2.1558600000000001
$ws1.Cells.Item(11, 4).Value = 0.1416
$ws1.Cells.Item(11, 5).Value = 0.70799999999999996
$ws1.Cells.Item(11, 6).Value = -3.5510470000000001
$ws1.Cells.Item(11, 7).Value = 29.378411
$ws1.Cells.Item(11, 8).Value = 3.05
$ws1.Cells.Item(11, 9).Value = 0.71
$ws1.Cells.Item(11, 10).Value = 1.28
$ws1.Cells.Item(11, 11).Value = 0

# Row 12 - MWUMBA
$ws1.Cells.Item(12, 1).Value = "MWUMBA"
$ws1.Cells.Item(12, 2).Value = 903
$ws1.Cells.Item(12, 3).Value = 3.6209999999999996
$ws1.Cells.Item(12, 4).Value = 0.18060000000000004
$ws1.Cells.Item(12, 5).Value = 0.90300000000000025
$ws1.Cells.Item(12, 6).Value = -3.2078739999999999
$ws1.Cells.Item(12, 7).Value = 29.872033999999999
$ws1.Cells.Item(12, 8).Value = 4.6100000000000003
$ws1.Cells.Item(12, 9).Value = 0.91
$ws1.Cells.Item(12, 10).Value = 1.39
$ws1.Cells.Item(12, 11).Value = 0

# Row 13 - NDAGO
$ws1.Cells.Item(13, 1).Value = "NDAGO"
$ws1.Cells.Item(13, 2).Value = 635
$ws1.Cells.Item(13, 3).Value = 0.95299999999999996
$ws1.Cells.Item(13, 4).Value = 0.127
$ws1.Cells.Item(13, 5).Value = 0.63500000000000001
$ws1.Cells.Item(13, 6).Value = -3.8171219999999999
$ws1.Cells.Item(13, 7).Value = 29.923649000000001
$ws1.Cells.Item(13, 8).Value = 3.47
$ws1.Cells.Item(13, 9).Value = 0.64
$ws1.Cells.Item(13, 10).Value = 1.1200000000000001
$ws1.Cells.Item(13, 11).Value = 0

# 4. Resize the "Tabelle1" table (and its autofilter) to the new extent
$lo = $ws1.ListObjects.Item(1)
$lo.Resize($ws1.Range("A1:K16"))

# 5. Update sheet selections / active sheet (Schools becomes the active tab)
$ws2.Range("D4").Select()
$ws1.Select()
$ws1.Range("B21").Select()
